$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data got re-sorted by Excels Sort feature (by column D) and a new
# row (TypeAccess) was inserted, so every row below 2 needs to be rewritten. Clear the
# whole data block first, then repopulate it with the final values.
$ws.Range("A2:L60").ClearContents()

$ws.Range("A1").Value = "Node"
$ws.Range("B1").Value = "Kadabra Node"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Extends"
$ws.Range("E1").Value = "Created"
$ws.Range("F1").Value = "Parsing"
$ws.Range("G1").Value = "Generating Code"
$ws.Range("H1").Value = "Notes"
$ws.Range("I1").Value = "JPModel (deprecated)"
$ws.Range("A2").Value = "EnumDecl"
$ws.Range("D2").Value = "ClassDecl"
$ws.Range("I2").Value = "x"
$ws.Range("L2").Value = "Partially based on:"
$ws.Range("A3").Value = "RecordDecl"
$ws.Range("D3").Value = "ClassDecl"
$ws.Range("I3").Value = "x"
$ws.Range("L3").Value = "https://docs.oracle.com/javase/specs/jls/se17/html/jls-19.html "
$ws.Range("A4").Value = "Import"
$ws.Range("D4").Value = "KadabraNode"
$ws.Range("I4").Value = "x"
$ws.Range("A5").Value = "Package"
$ws.Range("D5").Value = "KadabraNode"
$ws.Range("I5").Value = "x"
$ws.Range("A6").Value = "LambdaExpr"
$ws.Range("D6").Value = "Expression"
$ws.Range("I6").Value = "x"
$ws.Range("A7").Value = "LiteralExpr"
$ws.Range("D7").Value = "Expression"
$ws.Range("I7").Value = "x"
$ws.Range("A8").Value = "DeclReference"
$ws.Range("D8").Value = "Expression"
$ws.Range("I8").Value = "x"
$ws.Range("A9").Value = "ArrayAccess"
$ws.Range("D9").Value = "Expression"
$ws.Range("I9").Value = "x"
$ws.Range("A10").Value = "TypeAccess"
$ws.Range("B10").Value = "CtTypeAccess"
$ws.Range("D10").Value = "Expression"
$ws.Range("A11").Value = "MethodInvocation"
$ws.Range("B11").Value = "CtInvocationImpl"
$ws.Range("D11").Value = "Expression"
$ws.Range("I11").Value = "x"
$ws.Range("A12").Value = "MethodReference"
$ws.Range("D12").Value = "Expression"
$ws.Range("I12").Value = "x"
$ws.Range("A13").Value = "NewExpr"
$ws.Range("D13").Value = "Expression"
$ws.Range("I13").Value = "x"
$ws.Range("A14").Value = "Operator"
$ws.Range("D14").Value = "Expression"
$ws.Range("I14").Value = "x"
$ws.Range("A15").Value = "CastExpr"
$ws.Range("D15").Value = "Expression"
$ws.Range("I15").Value = "x"
$ws.Range("A16").Value = "SwitchExpr"
$ws.Range("D16").Value = "Expression"
$ws.Range("I16").Value = "x"
$ws.Range("A17").Value = "Decl"
$ws.Range("D17").Value = "KadabraNode"
$ws.Range("H17").Value = "Only elements that can be referenced should be Decl"
$ws.Range("I17").Value = "x"
$ws.Range("A18").Value = "Type"
$ws.Range("D18").Value = "KadabraNode"
$ws.Range("I18").Value = "x"
$ws.Range("A19").Value = "CompilationUnit"
$ws.Range("D19").Value = "KadabraNode"
$ws.Range("E19").Value = "x"
$ws.Range("I19").Value = "x"
$ws.Range("A20").Value = "Expression"
$ws.Range("D20").Value = "KadabraNode"
$ws.Range("I20").Value = "x"
$ws.Range("A21").Value = "App"
$ws.Range("D21").Value = "KadabraNode"
$ws.Range("E21").Value = "x"
$ws.Range("I21").Value = "x"
$ws.Range("A22").Value = "Statement"
$ws.Range("D22").Value = "KadabraNode"
$ws.Range("I22").Value = "x"
$ws.Range("A23").Value = "ClassLiteral"
$ws.Range("D23").Value = "LiteralExpr"
$ws.Range("I23").Value = "x"
$ws.Range("A24").Value = "ConstructorDecl"
$ws.Range("B24").Value = "CtConstructor"
$ws.Range("D24").Value = "MethodDecl"
$ws.Range("I24").Value = "x"
$ws.Range("A25").Value = "TypeDecl"
$ws.Range("B25").Value = "CtType"
$ws.Range("D25").Value = "Decl"
$ws.Range("I25").Value = "x"
$ws.Range("A26").Value = "MethodDecl"
$ws.Range("B26").Value = "CtMethod"
$ws.Range("D26").Value = "Decl"
$ws.Range("E26").Value = "x"
$ws.Range("I26").Value = "x"
$ws.Range("A27").Value = "FieldDecl"
$ws.Range("D27").Value = "Decl"
$ws.Range("I27").Value = "x"
$ws.Range("A28").Value = "VarDecl"
$ws.Range("D28").Value = "Decl"
$ws.Range("I28").Value = "x"
$ws.Range("A29").Value = "EnumConstantDecl"
$ws.Range("D29").Value = "Decl"
$ws.Range("I29").Value = "x"
$ws.Range("A30").Value = "UnaryOperator"
$ws.Range("D30").Value = "Operator"
$ws.Range("I30").Value = "x"
$ws.Range("A31").Value = "BinaryOperator"
$ws.Range("D31").Value = "Operator"
$ws.Range("I31").Value = "x"
$ws.Range("A32").Value = "TernaryOperator"
$ws.Range("D32").Value = "Operator"
$ws.Range("I32").Value = "x"
$ws.Range("A33").Value = "BlockStmt"
$ws.Range("B33").Value = "CtBlock"
$ws.Range("D33").Value = "Statement"
$ws.Range("E33").Value = "x"
$ws.Range("I33").Value = "x"
$ws.Range("A34").Value = "TryStmt"
$ws.Range("D34").Value = "Statement"
$ws.Range("I34").Value = "x"
$ws.Range("A35").Value = "ExprStmt"
$ws.Range("D35").Value = "Statement"
$ws.Range("I35").Value = "x"
$ws.Range("A36").Value = "DeclStmt"
$ws.Range("D36").Value = "Statement"
$ws.Range("I36").Value = "x"
$ws.Range("A37").Value = "LabelStmt"
$ws.Range("D37").Value = "Statement"
$ws.Range("I37").Value = "x"
$ws.Range("A38").Value = "IfStmt"
$ws.Range("D38").Value = "Statement"
$ws.Range("I38").Value = "x"
$ws.Range("A39").Value = "LoopStmt"
$ws.Range("D39").Value = "Statement"
$ws.Range("I39").Value = "x"
$ws.Range("A40").Value = "ForStmt"
$ws.Range("D40").Value = "LoopStmt"
$ws.Range("I40").Value = "x"
$ws.Range("A41").Value = "WhileStmt"
$ws.Range("D41").Value = "LoopStmt"
$ws.Range("I41").Value = "x"
$ws.Range("A42").Value = "EmptyStmt"
$ws.Range("D42").Value = "Statement"
$ws.Range("I42").Value = "x"
$ws.Range("A43").Value = "AssertStmt"
$ws.Range("D43").Value = "Statement"
$ws.Range("I43").Value = "x"
$ws.Range("A44").Value = "SwitchStmt"
$ws.Range("D44").Value = "Statement"
$ws.Range("I44").Value = "x"
$ws.Range("A45").Value = "DoStmt"
$ws.Range("D45").Value = "LoopStmt"
$ws.Range("I45").Value = "x"
$ws.Range("A46").Value = "BreakStmt"
$ws.Range("D46").Value = "Statement"
$ws.Range("I46").Value = "x"
$ws.Range("A47").Value = "ContinueStmt"
$ws.Range("D47").Value = "Statement"
$ws.Range("I47").Value = "x"
$ws.Range("A48").Value = "ReturnStmt"
$ws.Range("D48").Value = "Statement"
$ws.Range("I48").Value = "x"
$ws.Range("A49").Value = "SynchronizedStmt"
$ws.Range("D49").Value = "Statement"
$ws.Range("I49").Value = "x"
$ws.Range("A50").Value = "ThrowStmt"
$ws.Range("D50").Value = "Statement"
$ws.Range("I50").Value = "x"
$ws.Range("A51").Value = "YieldStmt"
$ws.Range("D51").Value = "Statement"
$ws.Range("I51").Value = "x"
$ws.Range("A52").Value = "PrimitiveType"
$ws.Range("D52").Value = "Type"
$ws.Range("I52").Value = "x"
$ws.Range("A53").Value = "ReferenceType"
$ws.Range("D53").Value = "Type"
$ws.Range("I53").Value = "x"
$ws.Range("A54").Value = "ClassDecl"
$ws.Range("B54").Value = "CtClass"
$ws.Range("D54").Value = "TypeDecl"
$ws.Range("E54").Value = "x"
$ws.Range("F54").Value = "x"
$ws.Range("H54").Value = "Can represent abstract class"
$ws.Range("I54").Value = "x"
$ws.Range("A55").Value = "InterfaceDecl"
$ws.Range("D55").Value = "TypeDecl"
$ws.Range("I55").Value = "x"
$ws.Range("A56").Value = "ParamDecl"
$ws.Range("D56").Value = "VarDecl"
$ws.Range("I56").Value = "x"
$ws.Range("A57").Value = "Comment"
$ws.Range("D57").Value = "KadabraNode"
$ws.Range("A58").Value = "CommentStmt"
$ws.Range("D58").Value = "Statement"
$ws.Range("A59").Value = "CaseStmt"
$ws.Range("D59").Value = "Statement"
$ws.Range("A60").Value = "Annotation"
$ws.Range("D60").Value = "Decl"

# Row count used by the "done" ratio grew by one (new TypeAccess row)
$ws.Range("K2").Formula = "=COUNTA(I:I)/COUNTA(A2:A56)"

# Restore the active selection used when the sheet was last saved
$ws.Range("F33").Select() | Out-Null
